$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3 - RandomForestRegressor
$ws.Range("B3").Value = 0.02107734011372254
$ws.Range("C3").Value = 0.02021425912045938
$ws.Range("D3").Value = 0.02065251971872556

# Row 4 - GradientBoostingRegressor
$ws.Range("B4").Value = 0.02403771557759567
$ws.Range("C4").Value = 0.0240350996315787
$ws.Range("D4").Value = 0.02403771557759567

# Row 5 - AdaBoostRegressor
$ws.Range("B5").Value = 0.1042252181250347
$ws.Range("C5").Value = 0.1188057510671667
$ws.Range("D5").Value = 0.1143316404626366
